$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: new value in S7 ---
$ws.Range("S7").Value = 4

# --- Row 21: new values in J21 and M21 ---
# J21 needs the same fill/border style as the other "grade" cells in that
# column (e.g. K6 / K7 use style index 13), so copy formatting first, then
# set the value.
$ws.Range("K6").Copy()
$ws.Range("J21").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("J21").Value = 5

# M21 just takes the plain column-M formatting (style 14), which is what a
# fresh value entry gets by default in this column.
$ws.Range("M21").Value = 3

# New S21 value
$ws.Range("S21").Value = 4

# --- Row 28: new values in S28 and T28 ---
$ws.Range("S28").Value = 4
$ws.Range("T28").Value = "авансорм"

# --- Sheet view: update frozen pane top-left cell and active selection ---
$sheetView = $ws.Application.ActiveWindow
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("E28").Select()
